$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPIEBP")
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2
